$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shared_otus")

# Row 3 - Unique OTUs: add percentages, convert to text
$ws.Range("B3").Value = "5 (2%)"
$ws.Range("C3").Value = "2 (1%)"
$ws.Range("D3").Value = "64 (15%)"
$ws.Range("E3").Value = "212 (58%)"

# Row 4 - Shared with Epiphytes: add percentages, convert to text
$ws.Range("C4").Value = "180 (89%)"
$ws.Range("D4").Value = "247 (59%)"
$ws.Range("E4").Value = "57 (16%)"

# Row 5 - Shared with Endophytes: add percentages, convert to text
$ws.Range("B5").Value = "180 (67%)"
$ws.Range("D5").Value = "187 (44%)"
$ws.Range("E5").Value = "43 (12%)"

# Row 6 - Shared with Litter: add percentages, convert to text
$ws.Range("B6").Value = "247 (92%)"
$ws.Range("C6").Value = "187 (93%)"
$ws.Range("E6").Value = "145 (40%)"

$wb.Save()
